$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- source row 28
$ws.Range("D2").Value = 44405
$ws.Range("H2").Value = "Madrigal"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 21000
$ws.Range("L2").Value = 22000
$ws.Range("M2").Value = 21500
$ws.Range("N2").Value = "`$/caja 40 unidades"
$ws.Range("O2").Value = "Región de Coquimbo"
$ws.Range("P2").Value = 538
$ws.Range("Q2").Value = 40

# Row 3 <- source row 22
$ws.Range("D3").Value = 44482
$ws.Range("H3").Value = "Madrigal"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14500
$ws.Range("N3").Value = "`$/caja 40 unidades"
$ws.Range("O3").Value = "Región de Coquimbo"
$ws.Range("P3").Value = 362
$ws.Range("Q3").Value = 40

# Row 4 <- source row 27
$ws.Range("D4").Value = 44769
$ws.Range("H4").Value = "Madrigal"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17500
$ws.Range("N4").Value = "`$/caja 40 unidades"
$ws.Range("O4").Value = "Región de Coquimbo"
$ws.Range("P4").Value = 438
$ws.Range("Q4").Value = 40

# Row 5 <- source row 26
$ws.Range("D5").Value = 44419
$ws.Range("H5").Value = "Symphony"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 21000
$ws.Range("L5").Value = 22000
$ws.Range("M5").Value = 21500
$ws.Range("N5").Value = "`$/caja 50 unidades"
$ws.Range("O5").Value = "Región de Coquimbo"
$ws.Range("P5").Value = 430
$ws.Range("Q5").Value = 50

# Row 6 <- source row 10
$ws.Range("D6").Value = 44468
$ws.Range("H6").Value = "Argentina(o)"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17500
$ws.Range("N6").Value = "`$/caja 50 unidades"
$ws.Range("O6").Value = "Región de Coquimbo"
$ws.Range("P6").Value = 350
$ws.Range("Q6").Value = 50

# Row 8 <- source row 16
$ws.Range("D8").Value = 44706
$ws.Range("H8").Value = "Madrigal"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 21000
$ws.Range("L8").Value = 22000
$ws.Range("M8").Value = 21500
$ws.Range("N8").Value = "`$/caja 40 unidades"
$ws.Range("O8").Value = "Región de Coquimbo"
$ws.Range("P8").Value = 538
$ws.Range("Q8").Value = 40

# Row 10 <- source row 3
$ws.Range("D10").Value = 44435
$ws.Range("H10").Value = "Madrigal"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 160
$ws.Range("K10").Value = 19000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 19500
$ws.Range("N10").Value = "`$/caja 40 unidades"
$ws.Range("O10").Value = "Región de Coquimbo"
$ws.Range("P10").Value = 488
$ws.Range("Q10").Value = 40

# Row 11 <- source row 25
$ws.Range("D11").Value = 44742
$ws.Range("H11").Value = "Madrigal"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 19000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 19500
$ws.Range("N11").Value = "`$/caja 40 unidades"
$ws.Range("O11").Value = "Región de Coquimbo"
$ws.Range("P11").Value = 488
$ws.Range("Q11").Value = 40

# Row 12 <- source row 6
$ws.Range("D12").Value = 44859
$ws.Range("H12").Value = "Madrigal"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 15500
$ws.Range("N12").Value = "`$/caja 40 unidades"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 388
$ws.Range("Q12").Value = 40

# Row 13 <- source row 17
$ws.Range("D13").Value = 44785
$ws.Range("H13").Value = "Argentina(o)"
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 160
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15500
$ws.Range("N13").Value = "`$/caja 50 unidades"
$ws.Range("O13").Value = "Región de Coquimbo"
$ws.Range("P13").Value = 310
$ws.Range("Q13").Value = 50

# Row 14 <- source row 30
$ws.Range("D14").Value = 44370
$ws.Range("H14").Value = "Argentina(o)"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 140
$ws.Range("K14").Value = 20000
$ws.Range("L14").Value = 21000
$ws.Range("M14").Value = 20429
$ws.Range("N14").Value = "`$/caja 50 unidades"
$ws.Range("O14").Value = "Región de Coquimbo"
$ws.Range("P14").Value = 409
$ws.Range("Q14").Value = 50

# Row 15 <- source row 31
$ws.Range("D15").Value = 44370
$ws.Range("H15").Value = "Madrigal"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 22000
$ws.Range("L15").Value = 23000
$ws.Range("M15").Value = 22500
$ws.Range("N15").Value = "`$/caja 40 unidades"
$ws.Range("O15").Value = "Región de Coquimbo"
$ws.Range("P15").Value = 562
$ws.Range("Q15").Value = 40

# Row 16 <- source row 12
$ws.Range("D16").Value = 44762
$ws.Range("H16").Value = "Madrigal"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 19000
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = 19500
$ws.Range("N16").Value = "`$/caja 40 unidades"
$ws.Range("O16").Value = "Región de Coquimbo"
$ws.Range("P16").Value = 488
$ws.Range("Q16").Value = 40

# Row 17 <- source row 29
$ws.Range("D17").Value = 44827
$ws.Range("H17").Value = "Madrigal"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 14000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 14500
$ws.Range("N17").Value = "`$/caja 40 unidades"
$ws.Range("O17").Value = "Región de Coquimbo"
$ws.Range("P17").Value = 362
$ws.Range("Q17").Value = 40

# Row 18 <- source row 23
$ws.Range("D18").Value = 44167
$ws.Range("H18").Value = "Española"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 160
$ws.Range("K18").Value = 13000
$ws.Range("L18").Value = 14000
$ws.Range("M18").Value = 13500
$ws.Range("N18").Value = "`$/caja 30 unidades"
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 450
$ws.Range("Q18").Value = 30

# Row 19 <- source row 4
$ws.Range("D19").Value = 44412
$ws.Range("H19").Value = "Symphony"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 240
$ws.Range("K19").Value = 21000
$ws.Range("L19").Value = 22000
$ws.Range("M19").Value = 21500
$ws.Range("N19").Value = "`$/caja 40 unidades"
$ws.Range("O19").Value = "Región de Coquimbo"
$ws.Range("P19").Value = 538
$ws.Range("Q19").Value = 40

# Row 20 <- source row 18
$ws.Range("D20").Value = 44806
$ws.Range("H20").Value = "Argentina(o)"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 250
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14500
$ws.Range("N20").Value = "`$/caja 40 unidades"
$ws.Range("O20").Value = "Provincia de Limarí"
$ws.Range("P20").Value = 362
$ws.Range("Q20").Value = 40

# Row 21 <- source row 19
$ws.Range("D21").Value = 44377
$ws.Range("H21").Value = "Madrigal"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 20000
$ws.Range("L21").Value = 21000
$ws.Range("M21").Value = 20333
$ws.Range("N21").Value = "`$/caja 40 unidades"
$ws.Range("O21").Value = "Región de Coquimbo"
$ws.Range("P21").Value = 508
$ws.Range("Q21").Value = 40

# Row 22 <- source row 20
$ws.Range("D22").Value = 44377
$ws.Range("H22").Value = "Symphony"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 21000
$ws.Range("L22").Value = 22000
$ws.Range("M22").Value = 21500
$ws.Range("N22").Value = "`$/caja 40 unidades"
$ws.Range("O22").Value = "Región de Coquimbo"
$ws.Range("P22").Value = 538
$ws.Range("Q22").Value = 40

# Row 23 <- source row 32
$ws.Range("D23").Value = 44356
$ws.Range("H23").Value = "Argentina(o)"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 120
$ws.Range("K23").Value = 19000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 19500
$ws.Range("N23").Value = "`$/caja 50 unidades"
$ws.Range("O23").Value = "Región de Coquimbo"
$ws.Range("P23").Value = 390
$ws.Range("Q23").Value = 50

# Row 24 <- source row 8
$ws.Range("D24").Value = 44398
$ws.Range("H24").Value = "Madrigal"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 170
$ws.Range("K24").Value = 21000
$ws.Range("L24").Value = 22000
$ws.Range("M24").Value = 21500
$ws.Range("N24").Value = "`$/caja 40 unidades"
$ws.Range("O24").Value = "Región de Coquimbo"
$ws.Range("P24").Value = 538
$ws.Range("Q24").Value = 40

# Row 25 <- source row 2
$ws.Range("D25").Value = 44483
$ws.Range("H25").Value = "Madrigal"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 120
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("N25").Value = "`$/caja 40 unidades"
$ws.Range("O25").Value = "Región de Coquimbo"
$ws.Range("P25").Value = 362
$ws.Range("Q25").Value = 40

# Row 26 <- source row 21
$ws.Range("D26").Value = 44489
$ws.Range("H26").Value = "Madrigal"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 13000
$ws.Range("L26").Value = 14000
$ws.Range("M26").Value = 13500
$ws.Range("N26").Value = "`$/caja 40 unidades"
$ws.Range("O26").Value = "Región de Coquimbo"
$ws.Range("P26").Value = 338
$ws.Range("Q26").Value = 40

# Row 27 <- source row 5
$ws.Range("D27").Value = 44160
$ws.Range("H27").Value = "Madrigal"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 160
$ws.Range("K27").Value = 14000
$ws.Range("L27").Value = 15000
$ws.Range("M27").Value = 14500
$ws.Range("N27").Value = "`$/caja 40 unidades"
$ws.Range("O27").Value = "Región de Coquimbo"
$ws.Range("P27").Value = 362
$ws.Range("Q27").Value = 40

# Row 28 <- source row 11
$ws.Range("D28").Value = 44433
$ws.Range("H28").Value = "Madrigal"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 160
$ws.Range("K28").Value = 19000
$ws.Range("L28").Value = 20000
$ws.Range("M28").Value = 19500
$ws.Range("N28").Value = "`$/caja 40 unidades"
$ws.Range("O28").Value = "Región de Coquimbo"
$ws.Range("P28").Value = 488
$ws.Range("Q28").Value = 40

# Row 29 <- source row 13
$ws.Range("D29").Value = 44384
$ws.Range("H29").Value = "Madrigal"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 80
$ws.Range("K29").Value = 21000
$ws.Range("L29").Value = 22000
$ws.Range("M29").Value = 21500
$ws.Range("N29").Value = "`$/caja 40 unidades"
$ws.Range("O29").Value = "Región de Coquimbo"
$ws.Range("P29").Value = 538
$ws.Range("Q29").Value = 40

# Row 30 <- source row 14
$ws.Range("D30").Value = 44384
$ws.Range("H30").Value = "Madrigal"
$ws.Range("I30").Value = "Segunda"
$ws.Range("J30").Value = 30
$ws.Range("K30").Value = 19000
$ws.Range("L30").Value = 20000
$ws.Range("M30").Value = 19333
$ws.Range("N30").Value = "`$/caja 50 unidades"
$ws.Range("O30").Value = "Región de Coquimbo"
$ws.Range("P30").Value = 387
$ws.Range("Q30").Value = 50

# Row 31 <- source row 15
$ws.Range("D31").Value = 44384
$ws.Range("H31").Value = "Symphony"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 20000
$ws.Range("L31").Value = 21000
$ws.Range("M31").Value = 20400
$ws.Range("N31").Value = "`$/caja 40 unidades"
$ws.Range("O31").Value = "Región de Coquimbo"
$ws.Range("P31").Value = 510
$ws.Range("Q31").Value = 40

# Row 32 <- source row 24
$ws.Range("D32").Value = 44426
$ws.Range("H32").Value = "Madrigal"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 150
$ws.Range("K32").Value = 19000
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = 19500
$ws.Range("N32").Value = "`$/caja 40 unidades"
$ws.Range("O32").Value = "Región de Coquimbo"
$ws.Range("P32").Value = 488
$ws.Range("Q32").Value = 40
